$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '25.992.35', '  -1.93%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.802.10', '  -2.04%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.004', '  +0.29%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '240.49', '  -7.58%  ')
    ,@(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.006', '  +0.58%  ')
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.5075', '  -3.35%  ')
    ,@(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2512', '  -21.38%  ')
    ,@(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06117', '  -9.90%  ')
    ,@(10, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.810.99', '  -1.37%  ')
    ,@(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06916', '  -10.83%  ')
    ,@(12, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '15.06', '  -19.69%  ')
    ,@(13, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6165', '  -21.44%  ')
    ,@(14, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '79.55', '  -9.47%  ')
    ,@(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.424', '  -11.75%  ')
    ,@(16, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.002', '  +0.13%  ')
    ,@(17, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.005', '  +0.40%  ')
    ,@(18, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '26.031.63', '  -1.83%  ')
    ,@(19, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '11.45', '  -17.29%  ')
    ,@(20, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.057.53', '  -0.75%  ')
    ,@(21, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000005919', '  -25.52%  ')
    ,@(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '3.987', '  -13.72%  ')
    ,@(23, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '5.231', '  -12.33%  ')
    ,@(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.017', '  -14.26%  ')
    ,@(25, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '130.73', '  -7.84%  ')
    ,@(26, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.890', '  -13.58%  ')
    ,@(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '14.64', '  -13.54%  ')
    ,@(28, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '99.54', '  -10.92%  ')
    ,@(29, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.166', '  -30.82%  ')
    ,@(30, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08247', '  -5.20%  ')
    ,@(31, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '3.650', '  -12.26%  ')
    ,@(32, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.754', '  -4.12%  ')
    ,@(33, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '1.005', '  +0.56%  ')
    ,@(34, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.204', '  -21.39%  ')
    ,@(35, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04332', '  -11.26%  ')
    ,@(36, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.050', '  -7.58%  ')
    ,@(37, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6335', '  -12.84%  ')
    ,@(38, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.902', '  -6.24%  ')
    ,@(39, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.104', '  -6.21%  ')
    ,@(40, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.006', '  +0.48%  ')
    ,@(41, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '99.94', '  -8.90%  ')
    ,@(42, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8039', '  -10.13%  ')
    ,@(43, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01459', '  -16.83%  ')
    ,@(44, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.3928', '  -17.92%  ')
    ,@(45, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '5.228', '  -11.96%  ')
    ,@(46, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.241', '  -18.75%  ')
    ,@(47, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05277', '  -9.76%  ')
    ,@(48, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '53.41', '  -10.29%  ')
    ,@(49, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1044', '  -15.32%  ')
    ,@(50, 'USDD', 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd', '1.008', '  +0.40%  ')
    ,@(51, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '29.43', '  -15.56%  ')
)

foreach ($row in $data) {
    $r = $row[0]

    $cB = $ws.Cells.Item($r, 2)
    $cC = $ws.Cells.Item($r, 3)
    $cD = $ws.Cells.Item($r, 4)
    $cE = $ws.Cells.Item($r, 5)

    # Force text number format so Excel does not reinterpret values such as
    # "1.004" or "25.992.35" as numbers/dates.
    $cB.NumberFormat = "@"
    $cC.NumberFormat = "@"
    $cD.NumberFormat = "@"
    $cE.NumberFormat = "@"

    $cB.Value = [string]$row[1]
    $cC.Value = [string]$row[2]
    $cD.Value = [string]$row[3]
    $cE.Value = [string]$row[4]

    # Restore the original (default) cell style now that the text value is set.
    $cB.Style = "Normal"
    $cC.Style = "Normal"
    $cD.Style = "Normal"
    $cE.Style = "Normal"
}
